$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix duplicated "类" typo in header labels (row 1)
$ws.Range("B1").Value = "其他服务类农村居民消费价格指数(上年=100)"
$ws.Range("D1").Value = "其他用品类农村居民消费价格指数(上年=100)"

# Add new row of data: 2021年 (copy formatting of the row above for the label cell)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 97.8
$ws.Range("C7").Value = 98.8
$ws.Range("D7").Value = 99.8

# Add new row of data: 2022年 (only the C column value is known)
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "2022年"
$ws.Range("C8").Value = 102
